$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: split the run "Dokumentacja" into three runs:
#           "Dokumentacj" + "s" + "a" (same visible text, same
#           formatting, but stored as three separate <w:r> elements).
# -----------------------------------------------------------------

$titlePara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "Dokumentacja*") {
        $titlePara = $para
    }
}

$titleStart = $titlePara.Range.Start

# Remove the first 11 characters ("Dokumentacj") from the run, leaving
# only the trailing "a" in the original paragraph (this paragraph keeps
# its original paragraph identity/properties).
$d.Range($titleStart, $titleStart + 11).Delete()

# Insert a new (empty) paragraph in front of it and fill it with
# "Dokumentacj".
$d.Range($titleStart, $titleStart).InsertParagraphBefore()
$titleIndex = $titlePara.Index
$d.Paragraphs($titleIndex).Range.InsertBefore("Dokumentacj")

# Insert another new paragraph between "Dokumentacj" and "a", filled
# with "s".
$midIndex = $titleIndex + 1
$midStart = $d.Paragraphs($midIndex).Range.Start
$d.Range($midStart, $midStart).InsertParagraphBefore()
$d.Paragraphs($midIndex).Range.InsertBefore("s")

# Merge "Dokumentacj" paragraph with "s" paragraph (delete the
# paragraph mark ending the "Dokumentacj" paragraph).
$mergeEnd1 = $d.Paragraphs($titleIndex).Range.End
$d.Range($mergeEnd1 - 1, $mergeEnd1).Delete()

# Merge the resulting "Dokumentacjs" paragraph with the trailing "a"
# paragraph (which is the original paragraph, so the final paragraph
# keeps the original paragraph's identity/properties).
$mergeEnd2 = $d.Paragraphs($titleIndex).Range.End
$d.Range($mergeEnd2 - 1, $mergeEnd2).Delete()

# -----------------------------------------------------------------
# Change 2: add a new bullet item after "Dodatkowo w aplikacji
#           została wykorzystana baza danych SQLite." reading
#           "Serwer aplikacji wykorzystuje Nginx".
# -----------------------------------------------------------------

$sqlitePara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "Dodatkowo w aplikacji*SQLite.*") {
        $sqlitePara = $para
    }
}

$sqlitePara.Range.InsertParagraphAfter()
$newParaIndex = $sqlitePara.Index + 1
$newPara = $d.Paragraphs($newParaIndex)
$newPara.Range.Text = "Serwer aplikacji wykorzystuje Nginx"
